$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells are kept as plain text so that values such as
# "576.56" or "6.76" are not auto-coerced into floating point numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.388.25'
$ws.Range('E2').Value = '  +4.79%  '

$ws.Range('D3').Value = '3.239.78'
$ws.Range('E3').Value = '  +2.55%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '576.56'
$ws.Range('E5').Value = '  +2.34%  '

$ws.Range('D6').Value = '179.11'
$ws.Range('E6').Value = '  +6.31%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  -2.46%  '

$ws.Range('D9').Value = '3.239.28'
$ws.Range('E9').Value = '  +2.71%  '

$ws.Range('E10').Value = '  +4.08%  '

$ws.Range('D11').Value = '6.76'

$ws.Range('D12').Value = '0.412'
$ws.Range('E12').Value = '  +4.70%  '

$ws.Range('D13').Value = '3.797.81'
$ws.Range('E13').Value = '  +2.54%  '

$ws.Range('E14').Value = '  +1.07%  '

$ws.Range('D15').Value = '27.86'
$ws.Range('E15').Value = '  +3.01%  '

$ws.Range('D16').Value = '67.345.90'
$ws.Range('E16').Value = '  +4.80%  '

$ws.Range('E17').Value = '  +2.61%  '

$ws.Range('D18').Value = '3.239.70'
$ws.Range('E18').Value = '  +2.83%  '

$ws.Range('E19').Value = '  +1.30%  '

$ws.Range('D20').Value = '13.34'
$ws.Range('E20').Value = '  +3.61%  '

$ws.Range('D21').Value = '373.63'
$ws.Range('E21').Value = '  +6.19%  '

$ws.Range('D22').Value = '7.58'
$ws.Range('E22').Value = '  +5.35%  '

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').Value = '71.10'
$ws.Range('E24').Value = '  +4.50%  '

$ws.Range('D25').Value = '0.508'
$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('E26').Value = '  +2.43%  '

$ws.Range('D27').Value = '9.63'
$ws.Range('E27').Value = '  +0.38%  '

$ws.Range('E28').Value = '  +3.16%  '

$ws.Range('E29').Value = '  +0.41%  '

$ws.Range('E30').Value = '  +4.39%  '

$ws.Range('E31').Value = '  +2.63%  '

$ws.Range('D32').Value = '22.59'
$ws.Range('E32').Value = '  +3.20%  '

$ws.Range('E34').Value = '  +6.41%  '

$ws.Range('D35').Value = '6.81'
$ws.Range('E35').Value = '  +3.34%  '

$ws.Range('D36').Value = '163.87'
$ws.Range('E36').Value = '  +6.58%  '

$ws.Range('E37').Value = '  +4.65%  '

$ws.Range('D38').Value = '0.861'
$ws.Range('E38').Value = '  +5.54%  '

$ws.Range('D39').Value = '1.85'
$ws.Range('E39').Value = '  +9.04%  '

$ws.Range('D40').Value = '6.86'
$ws.Range('E40').Value = '  +15.48%  '

$ws.Range('D41').Value = '26.73'
$ws.Range('E41').Value = '  +1.47%  '

$ws.Range('D42').Value = '362.49'
$ws.Range('E42').Value = '  +13.68%  '

$ws.Range('E43').Value = '  +5.15%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '4.40'
$ws.Range('E44').Value = '  +5.83%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.703.18'
$ws.Range('E45').Value = '  +3.04%  '

$ws.Range('D46').Value = '25.71'
$ws.Range('E46').Value = '  +8.30%  '

$ws.Range('D47').Value = '40.42'
$ws.Range('E47').Value = '  +2.76%  '

$ws.Range('E48').Value = '  +3.88%  '

$ws.Range('E50').Value = '  +0.58%  '

$ws.Range('D51').Value = '0.996'
$ws.Range('E51').Value = '  +6.27%  '
